{"js": "const replacements = [\n  [\"2024-11-10 Sunday\", \"2024-11-11 Monday\"],\n  [\"39\u00f73=13, 0\", \"87\u00f75=17, 2\"],\n  [\"48\u00f75=9, 3\", \"61\u00f72=30, 1\"],\n  [\"41\u00f75=8, 1\", \"22\u00f76=3, 4\"],\n  [\"30\u00f76=5, 0\", \"54\u00f76=9, 0\"],\n  [\"80\u00f76=13, 2\", \"69\u00f75=13, 4\"],\n  [\"42\u00f74=10, 2\", \"15\u00f77=2, 1\"],\n  [\"55\u00f76=9, 1\", \"76\u00f75=15, 1\"],\n  [\"20\u00f77=2, 6\", \"65\u00f79=7, 2\"],\n  [\"47\u00f72=23, 1\", \"68\u00f74=17, 0\"],\n  [\"82\u00f78=10, 2\", \"41\u00f76=6, 5\"],\n  [\"90\u00f79=10, 0\", \"38\u00f77=5, 3\"],\n  [\"42\u00f72=21, 0\", \"96\u00f78=12, 0\"],\n  [\"53\u00f78=6, 5\", \"48\u00f75=9, 3\"],\n  [\"55\u00f74=13, 3\", \"13\u00f77=1, 6\"],\n  [\"22\u00f77=3, 1\", \"16\u00f76=2, 4\"],\n  [\"83\u00f76=13, 5\", \"44\u00f76=7, 2\"],\n  [\"46\u00f77=6, 4\", \"37\u00f77=5, 2\"],\n  [\"75\u00f76=12, 3\", \"16\u00f75=3, 1\"],\n  [\"42\u00f73=14, 0\", \"37\u00f76=6, 1\"],\n  [\"24\u00f78=3, 0\", \"10\u00f72=5, 0\"],\n  [\"42\u00f75=8, 2\", \"16\u00f78=2, 0\"],\n  [\"44\u00f79=4, 8\", \"33\u00f73=11, 0\"],\n  [\"86\u00f72=43, 0\", \"65\u00f75=13, 0\"],\n  [\"33\u00f75=6, 3\", \"48\u00f78=6, 0\"],\n  [\"26\u00f77=3, 5\", \"23\u00f74=5, 3\"],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text([string]$find, [string]$replace) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $find\n    $r.Find.Replacement.Text = $replace\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 0\n    $r.Find.MatchCase = $true\n    $r.Find.MatchWholeWord = $false\n    $r.Find.Execute([ref]$find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\nReplace-Text \"2024-11-10 Sunday\" \"2024-11-11 Monday\"\nReplace-Text \"39\u00f73=13, 0\" \"87\u00f75=17, 2\"\nReplace-Text \"48\u00f75=9, 3\" \"61\u00f72=30, 1\"\nReplace-Text \"41\u00f75=8, 1\" \"22\u00f76=3, 4\"\nReplace-Text \"30\u00f76=5, 0\" \"54\u00f76=9, 0\"\nReplace-Text \"80\u00f76=13, 2\" \"69\u00f75=13, 4\"\nReplace-Text \"42\u00f74=10, 2\" \"15\u00f77=2, 1\"\nReplace-Text \"55\u00f76=9, 1\" \"76\u00f75=15, 1\"\nReplace-Text \"20\u00f77=2, 6\" \"65\u00f79=7, 2\"\nReplace-Text \"47\u00f72=23, 1\" \"68\u00f74=17, 0\"\nReplace-Text \"82\u00f78=10, 2\" \"41\u00f76=6, 5\"\nReplace-Text \"90\u00f79=10, 0\" \"38\u00f77=5, 3\"\nReplace-Text \"42\u00f72=21, 0\" \"96\u00f78=12, 0\"\nReplace-Text \"53\u00f78=6, 5\" \"48\u00f75=9, 3\"\nReplace-Text \"55\u00f74=13, 3\" \"13\u00f77=1, 6\"\nReplace-Text \"22\u00f77=3, 1\" \"16\u00f76=2, 4\"\nReplace-Text \"83\u00f76=13, 5\" \"44\u00f76=7, 2\"\nReplace-Text \"46\u00f77=6, 4\" \"37\u00f77=5, 2\"\nReplace-Text \"75\u00f76=12, 3\" \"16\u00f75=3, 1\"\nReplace-Text \"42\u00f73=14, 0\" \"37\u00f76=6, 1\"\nReplace-Text \"24\u00f78=3, 0\" \"10\u00f72=5, 0\"\nReplace-Text \"42\u00f75=8, 2\" \"16\u00f78=2, 0\"\nReplace-Text \"44\u00f79=4, 8\" \"33\u00f73=11, 0\"\nReplace-Text \"86\u00f72=43, 0\" \"65\u00f75=13, 0\"\nReplace-Text \"33\u00f75=6, 3\" \"48\u00f78=6, 0\"\nReplace-Text \"26\u00f77=3, 5\" \"23\u00f74=5, 3\"\n"}
